$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '29.326.19'
$ws.Range('D2').ClearFormats()
$ws.Range('E2').Value = '  -0.18%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '1.843.76'
$ws.Range('D3').ClearFormats()
$ws.Range('E3').Value = '  -0.30%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '0.9989'
$ws.Range('D4').ClearFormats()
$ws.Range('E4').Value = '  -0.08%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '239.96'
$ws.Range('D5').ClearFormats()
$ws.Range('E5').Value = '  -0.24%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.6266'
$ws.Range('D6').ClearFormats()
$ws.Range('E6').Value = '  -0.11%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.9998'
$ws.Range('D7').ClearFormats()
$ws.Range('E7').Value = '  -0.02%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.07453'
$ws.Range('D8').ClearFormats()
$ws.Range('E8').Value = '  -2.06%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.2895'
$ws.Range('D9').ClearFormats()
$ws.Range('E9').Value = '  -0.32%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '24.42'
$ws.Range('D10').ClearFormats()
$ws.Range('E10').Value = '  -1.12%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.07738'
$ws.Range('D11').ClearFormats()
$ws.Range('E11').Value = '  -0.13%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '1.843.86'
$ws.Range('D12').ClearFormats()
$ws.Range('E12').Value = '  -2.35%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '4.983'
$ws.Range('D13').ClearFormats()
$ws.Range('E13').Value = '  -0.88%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '0.6787'
$ws.Range('D14').ClearFormats()
$ws.Range('E14').Value = '  -0.04%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.00001044'
$ws.Range('D15').ClearFormats()
$ws.Range('E15').Value = '  -1.10%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '81.85'
$ws.Range('D16').ClearFormats()
$ws.Range('E16').Value = '  -1.46%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '6.189'
$ws.Range('D17').ClearFormats()
$ws.Range('E17').Value = '  +0.84%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '29.371.17'
$ws.Range('D18').ClearFormats()
$ws.Range('E18').Value = '  -0.11%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '228.01'
$ws.Range('D19').ClearFormats()
$ws.Range('E19').Value = '  -0.48%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '12.29'
$ws.Range('D20').ClearFormats()
$ws.Range('E20').Value = '  -0.54%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '0.9995'
$ws.Range('D21').ClearFormats()
$ws.Range('E21').Value = '  -0.03%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '7.505'
$ws.Range('D22').ClearFormats()
$ws.Range('E22').Value = '  +0.24%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '1.000'
$ws.Range('D23').ClearFormats()
$ws.Range('E23').Value = '  -0.04%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '159.03'
$ws.Range('D24').ClearFormats()
$ws.Range('E24').Value = '  +0.10%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '8.447'
$ws.Range('D25').ClearFormats()
$ws.Range('E25').Value = '  +0.09%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '0.1364'
$ws.Range('D26').ClearFormats()
$ws.Range('E26').Value = '  -1.59%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '17.48'
$ws.Range('D27').ClearFormats()
$ws.Range('E27').Value = '  -1.25%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '0.06469'
$ws.Range('D28').ClearFormats()
$ws.Range('E28').Value = '  +15.35%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '1.420'
$ws.Range('D29').ClearFormats()
$ws.Range('E29').Value = '  -1.01%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '1.483'
$ws.Range('D30').ClearFormats()
$ws.Range('E30').Value = '  +1.07%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '4.082'
$ws.Range('D31').ClearFormats()
$ws.Range('E31').Value = '  -0.54%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '4.082'
$ws.Range('D32').ClearFormats()
$ws.Range('E32').Value = '  +0.29%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '1.827'
$ws.Range('D33').ClearFormats()
$ws.Range('E33').Value = '  -0.26%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '1.139'
$ws.Range('D34').ClearFormats()
$ws.Range('E34').Value = '  -1.99%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.6937'
$ws.Range('D35').ClearFormats()
$ws.Range('E35').Value = '  +0.17%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '2.581'
$ws.Range('D36').ClearFormats()
$ws.Range('E36').Value = '  -0.16%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '1.261.41'
$ws.Range('D37').ClearFormats()
$ws.Range('E37').Value = '  +2.02%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '2.838'
$ws.Range('D38').ClearFormats()
$ws.Range('E38').Value = '  +4.00%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.01833'
$ws.Range('D39').ClearFormats()
$ws.Range('E39').Value = '  +1.91%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '6.780'
$ws.Range('D40').ClearFormats()
$ws.Range('E40').Value = '  +6.35%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.9178'
$ws.Range('D41').ClearFormats()
$ws.Range('E41').Value = '  +1.99%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.9987'
$ws.Range('D42').ClearFormats()
$ws.Range('E42').Value = '  -0.12%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '2.004.56'
$ws.Range('D43').ClearFormats()
$ws.Range('E43').Value = '  +1.19%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '101.31'
$ws.Range('D44').ClearFormats()
$ws.Range('E44').Value = '  -0.12%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '66.02'
$ws.Range('D45').ClearFormats()
$ws.Range('E45').Value = '  +0.82%  '
$ws.Range('B46').Value = 'BabyDogeCoin'
$ws.Range('C46').Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.00000000119'
$ws.Range('D46').ClearFormats()
$ws.Range('E46').Value = '  +0.70%  '
$ws.Range('B47').Value = 'RenderToken'
$ws.Range('C47').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '1.732'
$ws.Range('D47').ClearFormats()
$ws.Range('E47').Value = '  +2.73%  '
$ws.Range('B48').Value = 'Aptos'
$ws.Range('C48').Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '7.054'
$ws.Range('D48').ClearFormats()
$ws.Range('E48').Value = '  -1.86%  '
$ws.Range('B49').Value = 'Algorand'
$ws.Range('C49').Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '0.1155'
$ws.Range('D49').ClearFormats()
$ws.Range('E49').Value = '  +1.15%  '
$ws.Range('B50').Value = 'EnergySwap'
$ws.Range('C50').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '9.008'
$ws.Range('D50').ClearFormats()
$ws.Range('E50').Value = '  +0.20%  '
$ws.Range('B51').Value = 'TheSandbox'
$ws.Range('C51').Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.3942'
$ws.Range('D51').ClearFormats()
$ws.Range('E51').Value = '  -1.28%  '
